$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix to include T_Can and T_Top every update
# Correct previously miscalculated CO2_Air_calc values
$ws.Range("D4").Value = 664.8000000000001
$ws.Range("D5").Value = 498.6
$ws.Range("D6").Value = 664.8000000000001

# Append new rows 7-11 with the same data format as existing rows
$newRows = @(
    @(5, "8/26/2018 19:20", 774, 831.0000000000001),
    @(6, "8/26/2018 19:25", 776, 831.0000000000001),
    @(7, "8/26/2018 19:30", 769, 831.0000000000001),
    @(8, "8/26/2018 19:35", 754, 831.0000000000001),
    @(9, "8/26/2018 19:50", 770, 831.0000000000001)
)

# Copy the style of the last existing data row's A cell so new index cells match formatting
$ws.Range("A6").Copy()

$rowIndex = 7
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $rowIndex++
}
